$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Where''s my order number 2028?'
$ws.Range("C2").Value = '{''orders'': [{''Product ID'': 16734, ''ProductName'': ''hewlett packard enterprise intel xeon l5530 2.40ghz 2.4ghz 8mb l3 processor'', ''Category'': ''CPUs'', ''Category ID'': 2615, ''OrderID'': 2028, ''CustomerID'': 969, ''OrderStatus'': ''Delivered'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-24 12:24:37.765490''}]}'

$ws.Range("B3").Value = 'Where is my order 1575?'
$ws.Range("C3").Value = '{''orders'': [{''Product ID'': 26041, ''ProductName'': ''samsung ms28j5215as solo microwave silver'', ''Category'': ''Microwaves'', ''Category ID'': 2618, ''OrderID'': 1575, ''CustomerID'': 1072, ''OrderStatus'': ''Cancelled'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-09-29 12:24:37.765490''}]}'

$ws.Range("B4").Value = 'Can I track my order number 2425?'
$ws.Range("C4").Value = '{''orders'': [{''Product ID'': 30776, ''ProductName'': ''bosch serie 4 wan28100gb 7kg 1400rpm freestanding washing machine white'', ''Category'': ''Washing Machines'', ''Category ID'': 2620, ''OrderID'': 2425, ''CustomerID'': 827, ''OrderStatus'': ''Pending'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-20 12:24:37.765490''}]}'

$ws.Range("B5").Value = 'What''s the status of order number 2659?'
$ws.Range("C5").Value = '{''orders'': [{''Product ID'': 16813, ''ProductName'': ''hp amd opteron quad core 2354 2.2ghz fio kit 2.2ghz 2mb l3 processor'', ''Category'': ''CPUs'', ''Category ID'': 2615, ''OrderID'': 2659, ''CustomerID'': 2326, ''OrderStatus'': ''Shipped'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-11 12:24:37.765490''}]}'

$ws.Range("B6").Value = 'When will my order 3003 arrive?'
$ws.Range("C6").Value = '{''orders'': [{''Product ID'': 37955, ''ProductName'': ''samsung brb260000ww built in fridge freezer 70/30 frost free'', ''Category'': ''Fridge Freezers'', ''Category ID'': 2622, ''OrderID'': 3003, ''CustomerID'': 297, ''OrderStatus'': ''Cancelled'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-09 12:24:37.765490''}]}'

$ws.Range("B7").Value = 'When is order 2885 going to be delivered?'
$ws.Range("C7").Value = '{''orders'': [{''Product ID'': 45206, ''ProductName'': ''swan sr11035bln freestanding retro under counter fridge freezer a rated 90 litre blue'', ''Category'': ''Fridges'', ''Category ID'': 2623, ''OrderID'': 2885, ''CustomerID'': 549, ''OrderStatus'': ''Shipped'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-02 12:24:37.765490''}]}'

$ws.Range("B8").Value = 'What''s the delivery date for order 2247?'
$ws.Range("C8").Value = '{''orders'': [{''Product ID'': 10726, ''ProductName'': ''lg 55uk6300pla 55 ultra hd 4k tv'', ''Category'': ''TVs'', ''Category ID'': 2614, ''OrderID'': 2247, ''CustomerID'': 563, ''OrderStatus'': ''Delivered'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-21 12:24:37.765490''}]}'

$ws.Range("B9").Value = 'How can I track order number 1847?'
$ws.Range("C9").Value = '{''orders'': [{''Product ID'': 22421, ''ProductName'': ''sony alpha 7 geh use'', ''Category'': ''Digital Cameras'', ''Category ID'': 2617, ''OrderID'': 1847, ''CustomerID'': 1452, ''OrderStatus'': ''Shipped'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-17 12:24:37.765490''}]}'

$ws.Range("B10").Value = 'Where is order 2846 now?'
$ws.Range("C10").Value = '{''orders'': [{''Product ID'': 17249, ''ProductName'': ''amd phenom ii x4 840 3.20ghz 2mb cache am3 quad core socket am3 retail box cpu processor inc fan and heatsink'', ''Category'': ''CPUs'', ''Category ID'': 2615, ''OrderID'': 2846, ''CustomerID'': 1086, ''OrderStatus'': ''Delivered'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-16 12:24:37.765490''}]}'

$ws.Range("B11").Value = 'I want to check the status of order 3273.'
$ws.Range("C11").Value = '{''orders'': [{''Product ID'': 24163, ''ProductName'': ''sharp r272wm compact microwave oven in white 800w 20 litre'', ''Category'': ''Microwaves'', ''Category ID'': 2618, ''OrderID'': 3273, ''CustomerID'': 1748, ''OrderStatus'': ''Cancelled'', ''ReturnEligible'': False, ''ShippingDate'': ''2024-10-15 12:24:37.765490''}]}'
